$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record was added to the daily log: insert a fresh row at 337,
# pushing the existing rows 337:402 down to 338:403 (dimension grows to R403).
$ws.Rows.Item(337).Insert()

# Fill the newly inserted row with the new "Camote" record (Perú origin).
$ws.Cells.Item(337, 1).Value  = 10
$ws.Cells.Item(337, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(337, 3).Value  = "La Araucanía"
$ws.Cells.Item(337, 4).Value  = 44504
$ws.Cells.Item(337, 5).Value  = 9
$ws.Cells.Item(337, 6).Value  = 100112045
$ws.Cells.Item(337, 7).Value  = "Zapallo"
$ws.Cells.Item(337, 8).Value  = "Camote"
$ws.Cells.Item(337, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(337, 10).Value = 650
$ws.Cells.Item(337, 11).Value = 800
$ws.Cells.Item(337, 12).Value = 800
$ws.Cells.Item(337, 13).Value = 800
$ws.Cells.Item(337, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(337, 15).Value = "Perú"
$ws.Cells.Item(337, 16).Value = 800
$ws.Cells.Item(337, 17).Value = 1
$ws.Cells.Item(337, 18).Value = "Hortaliza"
